# Apply weekly update: insert a new record row at row 23 (shifting the
# existing rows 23:93 down to 24:94), then populate the newly inserted
# row with the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 23; this pushes the old
# row 23 (and everything below it, through row 93) down to row 24..94,
# extending the used range to row 94, matching the new dimension A1:R94.
$ws.Rows("23:23").Insert()

# Populate the newly inserted row 23 with the new weekly record.
$ws.Range("A23").Value = 4
$ws.Range("B23").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C23").Value = "Los Lagos"
$ws.Range("D23").Value = 44414
$ws.Range("D23").NumberFormat = $ws.Range("D24").NumberFormat
$ws.Range("E23").Value = 10
$ws.Range("F23").Value = 100112039
$ws.Range("G23").Value = "Ciboulette"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 240
$ws.Range("K23").Value = 4500
$ws.Range("L23").Value = 4500
$ws.Range("M23").Value = 4500
$ws.Range("N23").Value = "`$/docena de atados"
$ws.Range("O23").Value = "Región Metropolitana"
$ws.Range("P23").Value = 1500
$ws.Range("Q23").Value = 3
$ws.Range("R23").Value = "Hortaliza"
